# Updated symbol list on Sun Feb 12 22:32:52 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values for the cryptos sheet.
# Target cells are stored as text (inline strings) in the original workbook,
# so we force the Text number format before assigning, then restore the
# default style so no stray formatting differences remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$editRange = $ws.Range("D2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "317.05"
$ws.Range("D3").Value = "41.06"
$ws.Range("E3").Value = "0.01%"
$ws.Range("D4").Value = "5.189"
$ws.Range("E4").Value = "1.33%"
$ws.Range("D5").Value = "0.07630"
$ws.Range("E5").Value = "-0.72%"
$ws.Range("D6").Value = "4.323"
$ws.Range("E6").Value = "1.18%"
$ws.Range("D7").Value = "1.660"
$ws.Range("E7").Value = "2.48%"
$ws.Range("D8").Value = "0.9337"
$ws.Range("E8").Value = "1.57%"
$ws.Range("D10").Value = "0.1242"
$ws.Range("E10").Value = "1.17%"
$ws.Range("D11").Value = "0.1827"
$ws.Range("E11").Value = "0.15%"
$ws.Range("D12").Value = "0.09058"
$ws.Range("E12").Value = "-0.89%"
$ws.Range("D13").Value = "0.04127"
$ws.Range("E13").Value = "-3.69%"
$ws.Range("E14").Value = "0.48%"
$ws.Range("D15").Value = "0.001273"
$ws.Range("E15").Value = "1.71%"
$ws.Range("D16").Value = "0.005921"
$ws.Range("E16").Value = "5.08%"
$ws.Range("D18").Value = "3.351"
$ws.Range("E18").Value = "-0.05%"
$ws.Range("D19").Value = "0.3363"
$ws.Range("E19").Value = "1.51%"
$ws.Range("D20").Value = "8.392"
$ws.Range("E20").Value = "21.26%"
$ws.Range("D21").Value = "0.1359"
$ws.Range("E21").Value = "-2.08%"
$ws.Range("D22").Value = "0.2875"
$ws.Range("E22").Value = "5.31%"
$ws.Range("D23").Value = "0.04032"
$ws.Range("E23").Value = "-0.29%"
$ws.Range("D24").Value = "0.001274"
$ws.Range("E24").Value = "0.62%"
$ws.Range("D25").Value = "0.004061"
$ws.Range("E25").Value = "-0.34%"
$ws.Range("D26").Value = "0.0001275"
$ws.Range("E26").Value = "0.58%"
$ws.Range("D38").Value = "0.02467"
$ws.Range("E38").Value = "0.25%"
$ws.Range("D39").Value = "0.05231"
$ws.Range("E39").Value = "-0.56%"
$ws.Range("D40").Value = "0.007784"
$ws.Range("E40").Value = "-0.61%"
$ws.Range("D41").Value = "0.1293"
$ws.Range("E41").Value = "-1.44%"
$ws.Range("D42").Value = "0.007084"
$ws.Range("E42").Value = "4.28%"
$ws.Range("D43").Value = "0.002102"
$ws.Range("E43").Value = "14.31%"
$ws.Range("D44").Value = "0.008229"
$ws.Range("E44").Value = "0.67%"
$ws.Range("D45").Value = "0.3432"
$ws.Range("E45").Value = "10.83%"
$ws.Range("D46").Value = "0.00006683"
$ws.Range("E46").Value = "-2.18%"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").Value = "0.78%"
$ws.Range("D48").Value = "0.1997"
$ws.Range("E48").Value = "-10.49%"
$ws.Range("D49").Value = "0.004225"
$ws.Range("E49").Value = "3.24%"
$ws.Range("D50").Value = "0.00002112"
$ws.Range("E50").Value = "0.78%"
$ws.Range("E51").Value = "0.78%"

$editRange.Style = "Normal"
